$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.712.53"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "3.117.01"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'244.84"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'616.99"
$ws.Range("D7").Value = "'1.10"
$ws.Range("E7").Value = "  -5.13%  "
$ws.Range("D8").Value = "'0.385"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "3.115.32"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'0.755"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "'5.62"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'34.78"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "91.496.33"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "3.703.58"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "3.111.19"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "'3.75"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "'5.80"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'449.98"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'0.0000202"
$ws.Range("E23").Value = "  -9.66%  "
$ws.Range("D24").Value = "'9.22"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "'5.86"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'89.51"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "'11.72"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "3.290.78"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  +17.96%  "
$ws.Range("D31").Value = "'0.226"
$ws.Range("E31").Value = "  -10.98%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").Value = "'0.167"
$ws.Range("E32").Value = "  -10.22%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'1.05"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").Value = "'0.176"
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("D35").Value = "'9.31"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'7.66"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "'26.25"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "'1.96"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "'3.94"
$ws.Range("E39").Value = "  -4.25%  "
$ws.Range("D40").Value = "'488.24"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'1.31"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "'0.435"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("D43").Value = "'3.43"
$ws.Range("E43").Value = "  -6.62%  "
$ws.Range("D44").Value = "'22.21"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'160.15"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").Value = "'0.698"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.43"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'44.16"
$ws.Range("E51").Value = "  -1.47%  "
